$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")
$ws2 = $wb.Worksheets.Item("ProductLoanOutput")

# Update the long product-name string used by both sheets' B1 cell (shared string reused in place)
$newProductName = "4219-RBI-EI-DB-DL-REC-RNI-INT-FFConMONonLASTSUN-FIFC-1-FFROP-DAILY-FIFR-1-MD-TR-1-OT-PE-1st"
$ws1.Range("B1").Value = $newProductName
$ws2.Range("B1").Value = $newProductName

# Change the shortname value on ProductLoanInput from numeric 4219 to text "421r"
$ws1.Range("B2").Value = "421r"

# Move the selection on ProductLoanInput to B12 (no longer multi-cell B2:B3)
$ws1.Range("B12").Select()

# Make ProductLoanOutput the active/selected tab
$ws2.Activate()
